$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: MSR100079 -> MSR100443
$ws.Range("A4").Value = "MSR100443"
$ws.Range("B4").Value = "60 DPD"
$ws.Range("C4").Value = "FNMA"
$ws.Range("D4").Value = 325958.04
$ws.Range("E4").Value = 488.13
$ws.Range("F4").Value = 326446.17
$ws.Range("G4").Value = 0.077
$ws.Range("H4").Value = 46018
$ws.Range("I4").Value = "Escrow Advance Capitalization"
$ws.Range("J4").Value = "AUTH-654369"

# Row 5: MSR100530 -> MSR100570
$ws.Range("A5").Value = "MSR100570"
$ws.Range("B5").Value = "90+ DPD"
$ws.Range("C5").Value = "FHLMC"
$ws.Range("D5").Value = 439738.41
$ws.Range("E5").Value = 399.3
$ws.Range("F5").Value = 440137.71
$ws.Range("G5").Value = 0.0746
$ws.Range("H5").Value = 46053
$ws.Range("I5").Value = "Deferred Interest Capitalization"
$ws.Range("J5").Value = "AUTH-253458"

# Row 6: MSR100578 -> MSR100726
$ws.Range("A6").Value = "MSR100726"
$ws.Range("B6").Value = "60 DPD"
$ws.Range("C6").Value = "FNMA"
$ws.Range("D6").Value = 321780.02
$ws.Range("E6").Value = 1056.05
$ws.Range("F6").Value = 322836.07
$ws.Range("G6").Value = 0.07630000000000001
$ws.Range("H6").Value = 46025
$ws.Range("I6").Value = "Escrow Advance Capitalization"
$ws.Range("J6").Value = "AUTH-980275"
